$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text before writing, so that
# numeric-looking strings (e.g. "310.50", "1.003") are not silently
# reinterpreted as numbers (which would drop trailing zeros / formatting).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.785.76'
$ws.Range("E2").Value = '  -2.61%  '
$ws.Range("D3").Value = '1.782.21'
$ws.Range("E3").Value = '  -2.18%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '310.50'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.5144'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '0.3790'
$ws.Range("E8").Value = '  -2.18%  '
$ws.Range("D9").Value = '0.07762'
$ws.Range("E9").Value = '  -8.06%  '
$ws.Range("D10").Value = '41.15'
$ws.Range("E10").Value = '  -1.64%  '
$ws.Range("D11").Value = '1.083'
$ws.Range("E11").Value = '  -2.45%  '
$ws.Range("D12").Value = '1.003'
$ws.Range("D13").Value = '6.196'
$ws.Range("E13").Value = '  -3.36%  '
$ws.Range("D14").Value = '20.10'
$ws.Range("E14").Value = '  -4.60%  '
$ws.Range("D15").Value = '1.779.36'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '7.155'
$ws.Range("E16").Value = '  -4.84%  '
$ws.Range("D17").Value = '91.43'
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("E18").Value = '  -6.09%  '
$ws.Range("D19").Value = '0.06544'
$ws.Range("E19").Value = '  -2.22%  '
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = '16.97'
$ws.Range("D22").Value = '5.903'
$ws.Range("E22").Value = '  -3.21%  '
$ws.Range("D23").Value = '27.839.17'
$ws.Range("E23").Value = '  -2.54%  '
$ws.Range("D24").Value = '10.99'
$ws.Range("E24").Value = '  -3.84%  '
$ws.Range("D25").Value = '2.240'
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").Value = '158.90'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '20.22'
$ws.Range("E27").Value = '  -4.32%  '
$ws.Range("D28").Value = '1.983.51'
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("D29").Value = '2.361'
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("D30").Value = '125.24'
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").Value = '0.1074'
$ws.Range("E31").Value = '  -0.82%  '
$ws.Range("D32").Value = '1.028'
$ws.Range("E32").Value = '  -6.32%  '
$ws.Range("D33").Value = '3.617'
$ws.Range("E33").Value = '  -1.59%  '
$ws.Range("D34").Value = '5.476'
$ws.Range("E34").Value = '  -4.83%  '
$ws.Range("D35").Value = '0.07099'
$ws.Range("E35").Value = '  -5.72%  '
$ws.Range("D36").Value = '0.02313'
$ws.Range("E36").Value = '  -2.28%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '8.680'
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = '0.2121'
$ws.Range("E38").Value = '  -4.87%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '11.55'
$ws.Range("E39").Value = '  +2.74%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = '5.008'
$ws.Range("E40").Value = '  -3.78%  '
$ws.Range("D41").Value = '0.6076'
$ws.Range("E41").Value = '  -4.10%  '
$ws.Range("D42").Value = '1.002'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").Value = '  -3.47%  '
$ws.Range("D44").Value = '1.320'
$ws.Range("E44").Value = '  -5.78%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.06'
$ws.Range("E45").Value = '  -3.61%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.5936'
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = '3.710'
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").Value = '127.78'
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("D49").Value = '1.215'
$ws.Range("E49").Value = '  +1.24%  '
$ws.Range("D50").Value = '1.893'
$ws.Range("E50").Value = '  -5.03%  '
$ws.Range("D51").Value = '0.06732'
$ws.Range("E51").Value = '  -3.66%  '

# Restore the default (unstyled) cell style now that the text is committed,
# matching the original workbook formatting (no explicit style on these cells).
$ws.Range("D2:E51").Style = "Normal"
